$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 13917014
$ws.Range("I6").Value = 27782778
$ws.Range("K6").Value = 83348334
$ws.Range("M6").Value = -83348222

$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("N12").ClearContents()

$ws.Range("H38").Value = 55555624
$ws.Range("I38").Value = 55555624
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 166666872
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = -166666500
$ws.Range("N38").ClearContents()

$ws.Range("H43").Value = 2408.6956
$ws.Range("I43").Value = 2613.2
$ws.Range("J43").Value = 2025.25
$ws.Range("K43").Value = 2613.2
$ws.Range("L43").Value = 2025.25
$ws.Range("M43").Value = -2544.2
$ws.Range("N43").Value = -2163.25

$ws.Range("H58").Value = 6495055
$ws.Range("I58").Value = 17857792
$ws.Range("J58").Value = 2062.1428
$ws.Range("K58").Value = 53573376
$ws.Range("L58").Value = 6186.428400000001
$ws.Range("M58").Value = -53573226
$ws.Range("N58").Value = -6486.428400000001

$ws.Range("H132").Value = 26389.146
$ws.Range("I132").Value = 4218.2915
$ws.Range("J132").Value = 79599.2
$ws.Range("K132").Value = 12654.8745
$ws.Range("L132").Value = 238797.6
$ws.Range("M132").Value = -10124.8745
$ws.Range("N132").Value = -243857.6

$ws.Range("H137").Value = 4280508.5
$ws.Range("I137").Value = 12821609
$ws.Range("K137").Value = 38464827
$ws.Range("M137").Value = -38462277

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1743.5454
$ws.Range("I110").Value = 1742.9
$ws.Range("K110").Value = 1742.9
$ws.Range("M110").Value = 302.0999999999999

$ws.Range("H132").Value = 12502581
$ws.Range("I132").Value = 20001754
$ws.Range("K132").Value = 60005262
$ws.Range("M132").Value = -60002732

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H109").Value = 20000
$ws.Range("J109").Value = 20000
$ws.Range("L109").Value = 20000
$ws.Range("N109").Value = -22774

$ws.Range("H115").Value = 19999.889
$ws.Range("J115").Value = 19999.889
$ws.Range("L115").Value = 19999.889
$ws.Range("N115").Value = -23133.889

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9531450
$ws.Range("I31").Value = 3459.8462
$ws.Range("J31").Value = 15161626
$ws.Range("K31").Value = 3459.8462
$ws.Range("L31").Value = 15161626
$ws.Range("M31").Value = -3164.8462
$ws.Range("N31").Value = -15162216

$ws.Range("H34").Value = 9531450
$ws.Range("I34").Value = 3459.8462
$ws.Range("J34").Value = 15161626
$ws.Range("K34").Value = 3459.8462
$ws.Range("L34").Value = 15161626
$ws.Range("M34").Value = -3257.8462
$ws.Range("N34").Value = -15162030

$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

$ws.Range("H132").Value = 638855.8
$ws.Range("I132").Value = 1877.0625
$ws.Range("J132").Value = 2337465.8
$ws.Range("K132").Value = 5631.1875
$ws.Range("L132").Value = 7012397.399999999
$ws.Range("M132").Value = -3101.1875
$ws.Range("N132").Value = -7017457.399999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 911.4375
$ws.Range("I68").Value = 976
$ws.Range("J68").Value = 882.0909
$ws.Range("K68").Value = 2928
$ws.Range("L68").Value = 2646.2727
$ws.Range("M68").Value = -2117
$ws.Range("N68").Value = -4268.2727

$ws.Range("H71").Value = 911.4375
$ws.Range("I71").Value = 976
$ws.Range("J71").Value = 882.0909
$ws.Range("K71").Value = 8784
$ws.Range("L71").Value = 7938.8181
$ws.Range("M71").Value = -4728
$ws.Range("N71").Value = -16050.8181

$ws.Range("H122").Value = 4143.8
$ws.Range("I122").Value = 505.33334
$ws.Range("J122").Value = 7782.2666
$ws.Range("K122").Value = 4548.00006
$ws.Range("L122").Value = 70040.39939999999
$ws.Range("M122").Value = -2098.00006
$ws.Range("N122").Value = -74940.39939999999

$ws.Range("H131").Value = 2228.5715
$ws.Range("J131").Value = 2720.625
$ws.Range("L131").Value = 8161.875
$ws.Range("N131").Value = -18241.875

$ws.Range("I132").Value = 1366.1666
$ws.Range("J132").Value = 2356.6
$ws.Range("K132").Value = 12295.4994
$ws.Range("L132").Value = 21209.4
$ws.Range("M132").Value = -9765.499400000001
$ws.Range("N132").Value = -26269.4

$ws.Range("H133").Value = 3757.5
$ws.Range("I133").Value = 3757.5
$ws.Range("K133").Value = 11272.5
$ws.Range("M133").Value = -6212.5

$ws.Range("H140").Value = 237773
$ws.Range("I140").Value = 445729
$ws.Range("J140").Value = 3822.5
$ws.Range("K140").Value = 1337187
$ws.Range("L140").Value = 11467.5
$ws.Range("M140").Value = -1332007
$ws.Range("N140").Value = -21827.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 15000
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents()

$ws.Range("H87").Value = 23833.334
$ws.Range("I87").Value = 18000
$ws.Range("K87").Value = 18000
$ws.Range("M87").Value = -16752

$ws.Range("H90").Value = 23833.334
$ws.Range("I90").Value = 18000
$ws.Range("K90").Value = 54000
$ws.Range("M90").Value = -47760

$ws.Range("H132").Value = 47622356
$ws.Range("I132").Value = 90911330
$ws.Range("J132").Value = 4490.7
$ws.Range("K132").Value = 272733990
$ws.Range("L132").Value = 13472.1
$ws.Range("M132").Value = -272731460
$ws.Range("N132").Value = -18532.1

$ws.Range("H134").Value = 29900
$ws.Range("J134").Value = 29900
$ws.Range("L134").Value = 89700
$ws.Range("N134").Value = -94770

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 5282.9697
$ws.Range("I46").Value = 854.63635
$ws.Range("J46").Value = 7497.136
$ws.Range("K46").Value = 854.63635
$ws.Range("L46").Value = 7497.136
$ws.Range("M46").Value = -666.63635
$ws.Range("N46").Value = -7873.136

$ws.Range("H93").Value = 1119.579
$ws.Range("I93").Value = 983
$ws.Range("J93").Value = 1182.6154
$ws.Range("K93").Value = 983
$ws.Range("L93").Value = 1182.6154
$ws.Range("M93").Value = 265
$ws.Range("N93").Value = -3678.6154

$ws.Range("H122").Value = 85725.25
$ws.Range("I122").Value = 102270.3
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 306810.9
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -304360.9
$ws.Range("N122").Value = -13900

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1839801.6
$ws.Range("I126").Value = 1839801.6
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 5519404.800000001
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -5516934.800000001
$ws.Range("N126").ClearContents()

$ws.Range("H128").Value = 48000
$ws.Range("J128").Value = 48000
$ws.Range("L128").Value = 48000
$ws.Range("N128").Value = -57960

$ws.Range("H129").Value = 40429
$ws.Range("J129").Value = 40429
$ws.Range("L129").Value = 40429
$ws.Range("N129").Value = -50429
